$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H41").Value = 1389
$ws.Range("I41").Value = 1000
$ws.Range("J41").Value = 1875.25
$ws.Range("K41").Value = 1000
$ws.Range("L41").Value = 1875.25
$ws.Range("M41").Value = -560
$ws.Range("N41").Value = -2755.25

$ws.Range("H64").Value = 3300.2222
$ws.Range("I64").Value = 3200.3333
$ws.Range("J64").Value = 3500
$ws.Range("K64").Value = 3200.3333
$ws.Range("L64").Value = 3500
$ws.Range("M64").Value = -2952.3333
$ws.Range("N64").Value = -3996

$ws.Range("H67").Value = 3300.2222
$ws.Range("I67").Value = 3200.3333
$ws.Range("J67").Value = 3500
$ws.Range("K67").Value = 3200.3333
$ws.Range("L67").Value = 3500
$ws.Range("M67").Value = -2342.3333
$ws.Range("N67").Value = -5216

$ws.Range("H74").Value = 7444.6665

$ws.Range("H77").Value = 7444.6665

$ws.Range("H86").Value = 15071.429
$ws.Range("I86").Value = 16166.667
$ws.Range("J86").Value = 8500
$ws.Range("K86").Value = 16166.667
$ws.Range("L86").Value = 8500
$ws.Range("M86").Value = -15043.667
$ws.Range("N86").Value = -10746

$ws.Range("H89").Value = 15071.429
$ws.Range("I89").Value = 16166.667
$ws.Range("J89").Value = 8500
$ws.Range("K89").Value = 80833.33499999999
$ws.Range("L89").Value = 42500
$ws.Range("M89").Value = -75217.33499999999
$ws.Range("N89").Value = -53732

$ws.Range("H92").Value = 896.9
$ws.Range("I92").Value = 719.5333000000001
$ws.Range("K92").Value = 719.5333000000001
$ws.Range("M92").Value = 528.4666999999999

$ws.Range("H101").Value = 219.66667
$ws.Range("I101").Value = 219.66667
$ws.Range("K101").Value = 659.00001
$ws.Range("M101").Value = 962.99999

$ws.Range("H103").Value = 1071
$ws.Range("I103").Value = 959.5
$ws.Range("J103").Value = 1108.1666
$ws.Range("K103").Value = 2878.5
$ws.Range("L103").Value = 3324.4998
$ws.Range("M103").Value = -2292.5
$ws.Range("N103").Value = -4496.4998

$ws.Range("H106").Value = 2593.44
$ws.Range("I106").Value = 1286.6666
$ws.Range("J106").Value = 3328.5
$ws.Range("K106").Value = 1286.6666
$ws.Range("L106").Value = 3328.5
$ws.Range("M106").Value = -655.6666
$ws.Range("N106").Value = -4590.5

$ws.Range("H107").Value = 45456428
$ws.Range("I107").Value = 55556620
$ws.Range("K107").Value = 55556620
$ws.Range("M107").Value = -55554700

$ws.Range("H115").Value = 19092256
$ws.Range("I115").Value = 22909708
$ws.Range("J115").Value = 5000
$ws.Range("K115").Value = 68729124
$ws.Range("L115").Value = 15000
$ws.Range("M115").Value = -68727557
$ws.Range("N115").Value = -18134

$ws.Range("H127").Value = 1389.7858
$ws.Range("I127").Value = 950.63635
$ws.Range("K127").Value = 2851.90905
$ws.Range("M127").Value = 2108.09095

$ws.Range("H129").Value = 1751.3889
$ws.Range("I129").Value = 1129.5834
$ws.Range("J129").Value = 2995
$ws.Range("K129").Value = 3388.7502
$ws.Range("L129").Value = 8985
$ws.Range("M129").Value = 1611.2498
$ws.Range("N129").Value = -18985

$ws.Range("H137").Value = 6867.684
$ws.Range("I137").Value = 6411.125
$ws.Range("K137").Value = 19233.375
$ws.Range("M137").Value = -16683.375

$ws.Range("H138").Value = 6141
$ws.Range("J138").Value = 6519.016
$ws.Range("L138").Value = 19557.048
$ws.Range("N138").Value = -29837.048


# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 46711.855
$ws.Range("I32").Value = 46932.84
$ws.Range("K32").Value = 46932.84
$ws.Range("M32").Value = -46645.84

$ws.Range("H110").Value = 7815003
$ws.Range("I110").Value = 11365459
$ws.Range("K110").Value = 11365459
$ws.Range("M110").Value = -11363414

$ws.Range("H132").Value = 10662.635
$ws.Range("I132").Value = 4468.564
$ws.Range("J132").Value = 29244.846
$ws.Range("K132").Value = 13405.692
$ws.Range("L132").Value = 87734.538
$ws.Range("M132").Value = -10875.692
$ws.Range("N132").Value = -92794.538


# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H105").Value = 30313252
$ws.Range("I105").Value = 43491224
$ws.Range("K105").Value = 43491224
$ws.Range("M105").Value = -43489477

$ws.Range("H140").Value = 100783.45
$ws.Range("J140").Value = 100783.45
$ws.Range("L140").Value = 100783.45
$ws.Range("N140").Value = -111143.45


# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H31").Value = 7556.7617
$ws.Range("I31").Value = 5415.857
$ws.Range("K31").Value = 5415.857
$ws.Range("M31").Value = -5120.857

$ws.Range("H34").Value = 7556.7617
$ws.Range("I34").Value = 5415.857
$ws.Range("K34").Value = 5415.857
$ws.Range("M34").Value = -5213.857

$ws.Range("H86").Value = 6483.952
$ws.Range("I86").Value = 5411.923
$ws.Range("K86").Value = 5411.923
$ws.Range("M86").Value = -4288.923

$ws.Range("H89").Value = 6483.952
$ws.Range("I89").Value = 5411.923
$ws.Range("K89").Value = 27059.615
$ws.Range("M89").Value = -21443.615

$ws.Range("H99").Value = 4341.8335
$ws.Range("I99").Value = 4134.9
$ws.Range("J99").Value = 4600.5
$ws.Range("K99").Value = 4134.9
$ws.Range("L99").Value = 4600.5
$ws.Range("M99").Value = -2636.9
$ws.Range("N99").Value = -7596.5

$ws.Range("H122").Value = 1065.5385
$ws.Range("I122").Value = 900.6667
$ws.Range("K122").Value = 2702.0001
$ws.Range("M122").Value = -252.0001000000002

$ws.Range("H126").Value = 4341.8335
$ws.Range("I126").Value = 4134.9
$ws.Range("J126").Value = 4600.5
$ws.Range("K126").Value = 12404.7
$ws.Range("L126").Value = 13801.5
$ws.Range("M126").Value = -9934.699999999999
$ws.Range("N126").Value = -18741.5

$ws.Range("H134").Value = 2206.1516
$ws.Range("I134").Value = 1623.5593
$ws.Range("J134").Value = 7116.5713
$ws.Range("K134").Value = 4870.6779
$ws.Range("L134").Value = 21349.7139
$ws.Range("M134").Value = -2335.6779
$ws.Range("N134").Value = -26419.7139

$ws.Range("H141").Value = 523889.66
$ws.Range("J141").Value = 542434.1
$ws.Range("L141").Value = 542434.1
$ws.Range("N141").Value = -552794.1


# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H113").Value = 3692.0557
$ws.Range("I113").Value = 2950
$ws.Range("K113").Value = 8850
$ws.Range("M113").Value = -6680

$ws.Range("H121").Value = 19610794
$ws.Range("J121").Value = 20835698
$ws.Range("L121").Value = 62507094
$ws.Range("N121").Value = -62509714

$ws.Range("H137").Value = 19528
$ws.Range("J137").Value = 3701.3333
$ws.Range("L137").Value = 11103.9999
$ws.Range("N137").Value = -21303.9999


# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H82").Value = 1500
$ws.Range("I82").Value = 1500
$ws.Range("J82").Value = 1500
$ws.Range("K82").Value = 1500
$ws.Range("L82").Value = 1500
$ws.Range("M82").Value = -1139
$ws.Range("N82").Value = -2222

$ws.Range("H85").Value = 1500
$ws.Range("I85").Value = 1500
$ws.Range("J85").Value = 1500
$ws.Range("K85").Value = 1500
$ws.Range("L85").Value = 1500
$ws.Range("M85").Value = -252
$ws.Range("N85").Value = -3996


# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H53").Value = 15000
$ws.Range("J53").Value = 15000
$ws.Range("L53").Value = 15000
$ws.Range("N53").Value = -16214

$ws.Range("H81").Value = 10154.393
$ws.Range("I81").Value = 7409.5386
$ws.Range("K81").Value = 14819.0772
$ws.Range("M81").Value = -13758.0772

$ws.Range("H84").Value = 10154.393
$ws.Range("I84").Value = 7409.5386
$ws.Range("K84").Value = 74095.386
$ws.Range("M84").Value = -68791.386

$ws.Range("H122").Value = 3596.9333
$ws.Range("I122").Value = 3353.8572
$ws.Range("K122").Value = 10061.5716
$ws.Range("M122").Value = -7611.571599999999

$ws.Range("H132").Value = 4228.4834
$ws.Range("I132").Value = 3694.2222
$ws.Range("J132").Value = 5029.875
$ws.Range("K132").Value = 11082.6666
$ws.Range("L132").Value = 15089.625
$ws.Range("M132").Value = -8552.6666
$ws.Range("N132").Value = -20149.625

